$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

# Locate the "DeCAF" column dynamically from the header row (row 1)
$decafCol = 8
for ($c = 1; $c -le $colCount; $c++) {
    $header = $ws.Cells.Item(1, $c).Value2
    if ($header -eq "DeCAF") {
        $decafCol = $c
    }
}

# Clear every cell in the DeCAF column whose value is "permCAF"
for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $decafCol)
    if ($cell.Value2 -eq "permCAF") {
        $cell.Value2 = ""
    }
}
